$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '89.686.13'
$ws.Range('E2').Value = '  +3.54%  '
$ws.Range('D3').Value = '3.197.77'
$ws.Range('E3').Value = '  +1.97%  '
$ws.Range('E4').Value = '  +0.22%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = '@'
$cell.Value = '217.30'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  +6.81%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = '@'
$cell.Value = '627.20'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  +3.45%  '
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.392'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  +7.26%  '
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.694'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  +6.38%  '
$ws.Range('E9').Value = '  +0.16%  '
$ws.Range('D10').Value = '3.193.29'
$ws.Range('E10').Value = '  +1.97%  '
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.577'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  +9.86%  '
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.179'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  +2.17%  '
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0000261'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  +9.34%  '
$ws.Range('B14').Value = 'Toncoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = '@'
$cell.Value = '5.42'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  +4.36%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = '@'
$cell.Value = '33.69'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  +6.03%  '
$ws.Range('D16').Value = '3.796.93'
$ws.Range('E16').Value = '  +2.22%  '
$ws.Range('D17').Value = '89.596.21'
$ws.Range('E17').Value = '  +3.83%  '
$ws.Range('D18').Value = '3.205.81'
$ws.Range('E18').Value = '  +2.72%  '
$ws.Range('B19').Value = 'SuiNetwork'
$ws.Range('C19').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = '@'
$cell.Value = '3.43'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  +16.71%  '
$ws.Range('B20').Value = 'PEPE'
$ws.Range('C20').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.0000227'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  +77.41%  '
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = '@'
$cell.Value = '13.55'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  +2.42%  '
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = '@'
$cell.Value = '437.69'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  +7.16%  '
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = '@'
$cell.Value = '8.67'
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  +3.59%  '
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = '@'
$cell.Value = '5.11'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  +1.37%  '
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = '@'
$cell.Value = '5.35'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  +5.52%  '
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = '@'
$cell.Value = '12.04'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  +3.78%  '
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = '@'
$cell.Value = '82.77'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  +14.09%  '
$ws.Range('D28').Value = '3.441.15'
$ws.Range('E28').Value = '  +4.16%  '
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  +0.08%  '
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.160'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  -0.09%  '
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.998'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  -0.15%  '
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = '@'
$cell.Value = '4.17'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  +41.43%  '
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = '@'
$cell.Value = '8.53'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  +4.25%  '
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = '@'
$cell.Value = '547.08'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  +3.20%  '
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = '@'
$cell.Value = '7.08'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  +9.67%  '
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.92'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  +4.78%  '
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.32'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  +3.31%  '
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = '@'
$cell.Value = '22.48'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  +4.84%  '
$ws.Range('E39').Value = '  +3.01%  '
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.129'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  -1.44%  '
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  +0.28%  '
$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.94'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  +3.56%  '
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.375'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  +2.78%  '
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = '@'
$cell.Value = '147.33'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  -0.99%  '
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = '@'
$cell.Value = '173.91'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  +2.27%  '
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = '@'
$cell.Value = '43.83'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  +2.47%  '
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.774'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  +13.69%  '
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.125'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  -0.74%  '
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = '@'
$cell.Value = '1.25'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  +1.57%  '
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = '@'
$cell.Value = '0.624'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  +8.02%  '
